# Adds a new "2022-Q4" sheet (fund-holdings detail) right after "总计",
# pushing all existing quarter sheets one slot later, and records the
# new quarter's summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by cloning the "2022-Q3" sheet (same
#    headers/column styling) and dropping it in right before "2022-Q3".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$templateIndex = $template.Index
$template.Copy($template)
$q4 = $wb.Worksheets.Item($templateIndex)
$q4.Name = "2022-Q4"

function Set-TextCell($cell, $text) {
    # Forces the cell to be stored as text (matches the source data, which
    # keeps these figures as strings) while keeping the default style.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# code, name, scale, position, ratio, marketValue, rank
$q4rows = @(
    @("011685", "创金合信先进装备股票A", "0.24", "80.29", "9.77", "0.0234", 3),
    @("011686", "创金合信先进装备股票C", "0.18", "80.29", "9.77", "0.0176", 3),
    @("004927", "中航军民融合精选混合C", "0.46", "58.03", "2.98", "0.0137", 8),
    @("004926", "中航军民融合精选混合A", "0.13", "58.03", "2.98", "0.0039", 8),
    @("004937", "中航混改精选混合C", "0.07", "74.18", "5.28", "0.0037", 9),
    @("004936", "中航混改精选混合A", "0.03", "74.18", "5.28", "0.0016", 9)
)

$templateLastRow = 6   # "2022-Q3" template had data in rows 2..6
$row = 2
foreach ($rec in $q4rows) {
    if ($row -gt $templateLastRow) {
        # Beyond the template's original extent: no row exists yet, so
        # clone column-A's formatting (bold/border index style) first.
        $q4.Cells.Item($templateLastRow, 1).Copy($q4.Cells.Item($row, 1))
    }
    $q4.Cells.Item($row, 1).Value = $row - 2

    Set-TextCell $q4.Cells.Item($row, 2) $rec[0]
    Set-TextCell $q4.Cells.Item($row, 3) $rec[1]
    Set-TextCell $q4.Cells.Item($row, 4) $rec[2]
    Set-TextCell $q4.Cells.Item($row, 5) $rec[3]
    Set-TextCell $q4.Cells.Item($row, 6) $rec[4]
    Set-TextCell $q4.Cells.Item($row, 7) $rec[5]
    $q4.Cells.Item($row, 8).Value = $rec[6]

    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2) Update "总计": insert the 2022-Q4 summary row at the top (row 2),
#    pushing every other quarter's B/C/D figures down by one row while
#    the A-column running index (0,1,2,...) is simply re-sequenced.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$lastDataRow = 9   # existing data occupies rows 2..9 before the edit
$oldB = @()
$oldC = @()
$oldD = @()
for ($r = 2; $r -le $lastDataRow; $r++) {
    $oldB += $total.Cells.Item($r, 2).Text
    $oldC += $total.Cells.Item($r, 3).Text
    $oldD += $total.Cells.Item($r, 4).Text
}

$newLastDataRow = $lastDataRow + 1   # one extra row appended at the bottom

# The very last row is brand new - clone column-A's index styling from the
# previous last row before writing into it.
$total.Cells.Item($lastDataRow, 1).Copy($total.Cells.Item($newLastDataRow, 1))

for ($r = 3; $r -le $newLastDataRow; $r++) {
    $idx = $r - 3
    $total.Cells.Item($r, 1).Value = $r - 2
    $total.Cells.Item($r, 2).Value = $oldB[$idx]
    $total.Cells.Item($r, 3).Value = [double]$oldC[$idx]
    $total.Cells.Item($r, 4).Value = [double]$oldD[$idx]
}

# Row 2 becomes the freshly-added 2022-Q4 summary (A2 is already 0).
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 6
$total.Cells.Item(2, 4).Value = 0.06
